$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Insured Name" column (D)
$ws.Columns("D:E").Insert()
# Insert two new columns after "Insured Name" (now column F), before "Quote No"
$ws.Columns("G:H").Insert()

# Fill in the header labels for the four newly inserted columns (set in
# this particular order so new shared-string entries are created in the
# same sequence as the authored workbook)
$ws.Range("E1").Value = "Class Description"
$ws.Range("H1").Value = "Referral Reason"
$ws.Range("D1").Value = "Class Code"
$ws.Range("G1").Value = "Industial Question Count"

# Match the new column widths (inputs chosen so the engine's internal
# pixel-quantization lands on the closest representable value to the
# target OOXML width: 21 and 25.28515625 respectively)
$ws.Columns("D:E").ColumnWidth = 20.15
$ws.Columns("G:H").ColumnWidth = 24.5

# Update the active selection as recorded in the saved workbook
$ws.Range("G8").Select()
